$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute(
    "neto iznos zarade od 450,00 eura za zaposlenog ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "neto iznos zarade od za zaposlenog ", 2
)
